# Auto-generated PowerShell/Excel COM-interop script
# Updates column F (想去人数 / "want to go" counts) values per the commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(3, 6).Value = 476
$ws.Cells.Item(4, 6).Value = 539
$ws.Cells.Item(5, 6).Value = 2636
$ws.Cells.Item(6, 6).Value = 20
$ws.Cells.Item(7, 6).Value = 100
$ws.Cells.Item(8, 6).Value = 106
$ws.Cells.Item(9, 6).Value = 1741
$ws.Cells.Item(10, 6).Value = 1741
$ws.Cells.Item(11, 6).Value = 1444
$ws.Cells.Item(12, 6).Value = 88
$ws.Cells.Item(13, 6).Value = 1478
$ws.Cells.Item(15, 6).Value = 42
$ws.Cells.Item(16, 6).Value = 1064
$ws.Cells.Item(17, 6).Value = 353
$ws.Cells.Item(18, 6).Value = 208
$ws.Cells.Item(19, 6).Value = 267
$ws.Cells.Item(20, 6).Value = 7716
$ws.Cells.Item(21, 6).Value = 8891
$ws.Cells.Item(23, 6).Value = 2
$ws.Cells.Item(24, 6).Value = 439
$ws.Cells.Item(26, 6).Value = 106
$ws.Cells.Item(32, 6).Value = 1557
$ws.Cells.Item(33, 6).Value = 37
$ws.Cells.Item(34, 6).Value = 280
$ws.Cells.Item(36, 6).Value = 34
$ws.Cells.Item(39, 6).Value = 838
$ws.Cells.Item(42, 6).Value = 384
$ws.Cells.Item(43, 6).Value = 288
$ws.Cells.Item(44, 6).Value = 239
$ws.Cells.Item(45, 6).Value = 101
$ws.Cells.Item(46, 6).Value = 239
$ws.Cells.Item(48, 6).Value = 214

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(20, 6).Value = 336

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 199
$ws.Cells.Item(3, 6).Value = 2681
$ws.Cells.Item(4, 6).Value = 317
$ws.Cells.Item(5, 6).Value = 169

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 199
$ws.Cells.Item(4, 6).Value = 476
$ws.Cells.Item(6, 6).Value = 317
$ws.Cells.Item(7, 6).Value = 169
$ws.Cells.Item(9, 6).Value = 2636
$ws.Cells.Item(10, 6).Value = 20
$ws.Cells.Item(11, 6).Value = 100
$ws.Cells.Item(12, 6).Value = 106
$ws.Cells.Item(13, 6).Value = 1741
$ws.Cells.Item(14, 6).Value = 1741
$ws.Cells.Item(15, 6).Value = 1444
$ws.Cells.Item(16, 6).Value = 88
$ws.Cells.Item(17, 6).Value = 1478
$ws.Cells.Item(19, 6).Value = 42
$ws.Cells.Item(22, 6).Value = 208
$ws.Cells.Item(23, 6).Value = 267
$ws.Cells.Item(24, 6).Value = 7716
$ws.Cells.Item(25, 6).Value = 8892
$ws.Cells.Item(27, 6).Value = 439
$ws.Cells.Item(28, 6).Value = 106
$ws.Cells.Item(31, 6).Value = 1557
$ws.Cells.Item(32, 6).Value = 37
$ws.Cells.Item(33, 6).Value = 280
$ws.Cells.Item(35, 6).Value = 34
$ws.Cells.Item(39, 6).Value = 838
$ws.Cells.Item(42, 6).Value = 384
$ws.Cells.Item(43, 6).Value = 288
$ws.Cells.Item(44, 6).Value = 239
$ws.Cells.Item(45, 6).Value = 101
$ws.Cells.Item(46, 6).Value = 239
$ws.Cells.Item(48, 6).Value = 214
$ws.Cells.Item(50, 6).Value = 336
